# Swap the order of "System" and "dnasr281@gmail.com" (and "admin@admin.com")
# within the "Recorded By" column (G) of the active worksheet.
#
# Rule observed in the target diff: any cell in column G whose text is
# exactly "System, dnasr281@gmail.com" becomes "dnasr281@gmail.com, System",
# and any cell whose text is exactly "admin@admin.com, dnasr281@gmail.com"
# becomes "dnasr281@gmail.com, admin@admin.com". All other combinations
# (single values, triples, or already-reordered pairs) are left untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$usedRange = $ws.UsedRange
$lastRow = $usedRange.Rows.Count

for ($row = 2; $row -le $lastRow; $row++) {
    $cell = $ws.Cells.Item($row, 7)   # Column G = "Recorded By"
    $value = $cell.Value2

    if ($value -eq "System, dnasr281@gmail.com") {
        $cell.Value = "dnasr281@gmail.com, System"
    }
    elseif ($value -eq "admin@admin.com, dnasr281@gmail.com") {
        $cell.Value = "dnasr281@gmail.com, admin@admin.com"
    }
}
